$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K quarterly data to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats / fonts / alignment from the (now shifted) old "D:E" data,
# which now lives in F:G, into the newly inserted D:E columns so they match the
# surrounding quarterly columns (date format for the header row, number format
# for the data rows, etc.)
$ws.Range("F7:G102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("F7").Value2 = 43281
$ws.Range("G7").Value2 = 43190
$ws.Range("H7").Value2 = 43100
$ws.Range("I7").Value2 = 43008
$ws.Range("J7").Value2 = 42916
$ws.Range("K7").Value2 = 42825
$ws.Range("L7").Value2 = 42735
$ws.Range("M7").Value2 = 42643
$ws.Range("D8").Value2 = 226300
$ws.Range("E8").Value2 = 223300
$ws.Range("F8").Value2 = 219500
$ws.Range("G8").Value2 = 212200
$ws.Range("H8").Value2 = 209200
$ws.Range("I8").Value2 = 208700
$ws.Range("J8").Value2 = 394100
$ws.Range("K8").Value2 = 194500
$ws.Range("L8").Value2 = 194600
$ws.Range("M8").Value2 = 192100
$ws.Range("D9").Value2 = 71700
$ws.Range("E9").Value2 = 73400
$ws.Range("F9").Value2 = 68700
$ws.Range("G9").Value2 = 67100
$ws.Range("H9").Value2 = 64900
$ws.Range("I9").Value2 = 68500
$ws.Range("J9").Value2 = 124600
$ws.Range("K9").Value2 = 60800
$ws.Range("L9").Value2 = 62400
$ws.Range("M9").Value2 = 62900
$ws.Range("D10").Value2 = 154600
$ws.Range("E10").Value2 = 149900
$ws.Range("F10").Value2 = 150800
$ws.Range("G10").Value2 = 145100
$ws.Range("H10").Value2 = 144300
$ws.Range("I10").Value2 = 140200
$ws.Range("J10").Value2 = 269500
$ws.Range("K10").Value2 = 133700
$ws.Range("L10").Value2 = 132200
$ws.Range("M10").Value2 = 129200
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"
$ws.Range("G12").Value2 = "NA"
$ws.Range("H12").Value2 = "NA"
$ws.Range("I12").Value2 = "NA"
$ws.Range("J12").Value2 = "NA"
$ws.Range("K12").Value2 = "NA"
$ws.Range("L12").Value2 = "NA"
$ws.Range("M12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 0
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("M13").Value2 = 0
$ws.Range("D14").Value2 = -1600
$ws.Range("E14").Value2 = -1300
$ws.Range("F14").Value2 = -1700
$ws.Range("G14").Value2 = -1500
$ws.Range("H14").Value2 = -1500
$ws.Range("I14").Value2 = -1100
$ws.Range("J14").Value2 = -3300
$ws.Range("K14").Value2 = -2200
$ws.Range("L14").Value2 = -2200
$ws.Range("M14").Value2 = -2300
$ws.Range("D15").Value2 = 89900
$ws.Range("E15").Value2 = 74100
$ws.Range("F15").Value2 = 73400
$ws.Range("G15").Value2 = 72500
$ws.Range("H15").Value2 = 70600
$ws.Range("I15").Value2 = 70000
$ws.Range("J15").Value2 = 136200
$ws.Range("K15").Value2 = 67400
$ws.Range("L15").Value2 = 67000
$ws.Range("M15").Value2 = 63800
$ws.Range("D17").Value2 = 170200
$ws.Range("E17").Value2 = 155600
$ws.Range("F17").Value2 = 149900
$ws.Range("G17").Value2 = 147600
$ws.Range("H17").Value2 = 143100
$ws.Range("I17").Value2 = 145800
$ws.Range("J17").Value2 = 276200
$ws.Range("K17").Value2 = 136200
$ws.Range("L17").Value2 = 136500
$ws.Range("M17").Value2 = 132500
$ws.Range("D18").Value2 = 56100
$ws.Range("E18").Value2 = 67700
$ws.Range("F18").Value2 = 69600
$ws.Range("G18").Value2 = 64600
$ws.Range("H18").Value2 = 66100
$ws.Range("I18").Value2 = 62900
$ws.Range("J18").Value2 = 117900
$ws.Range("K18").Value2 = 58300
$ws.Range("L18").Value2 = 58100
$ws.Range("M18").Value2 = 59600
$ws.Range("D20").Value2 = 900
$ws.Range("E20").Value2 = 1400
$ws.Range("F20").Value2 = 700
$ws.Range("G20").Value2 = 900
$ws.Range("H20").Value2 = 600
$ws.Range("I20").Value2 = 1100
$ws.Range("J20").Value2 = 1000
$ws.Range("K20").Value2 = 400
$ws.Range("L20").Value2 = 600
$ws.Range("M20").Value2 = 12600
$ws.Range("D21").Value2 = 147000
$ws.Range("E21").Value2 = 143200
$ws.Range("F21").Value2 = 143700
$ws.Range("G21").Value2 = 138000
$ws.Range("H21").Value2 = 137300
$ws.Range("I21").Value2 = 133900
$ws.Range("J21").Value2 = 255100
$ws.Range("K21").Value2 = 126100
$ws.Range("L21").Value2 = 125700
$ws.Range("M21").Value2 = 136100
$ws.Range("D22").Value2 = 33500
$ws.Range("E22").Value2 = 33700
$ws.Range("F22").Value2 = 33300
$ws.Range("G22").Value2 = 32900
$ws.Range("H22").Value2 = 34700
$ws.Range("I22").Value2 = 35500
$ws.Range("J22").Value2 = 75000
$ws.Range("K22").Value2 = 37000
$ws.Range("L22").Value2 = 36300
$ws.Range("M22").Value2 = 36500
$ws.Range("D23").Value2 = 23500
$ws.Range("E23").Value2 = 35400
$ws.Range("F23").Value2 = 37000
$ws.Range("G23").Value2 = 32600
$ws.Range("H23").Value2 = 32000
$ws.Range("I23").Value2 = 28500
$ws.Range("J23").Value2 = 43900
$ws.Range("K23").Value2 = 21800
$ws.Range("L23").Value2 = 22500
$ws.Range("M23").Value2 = 35800
$ws.Range("D24").Value2 = 0
$ws.Range("E24").Value2 = 0
$ws.Range("F24").Value2 = 0
$ws.Range("G24").Value2 = 0
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = 0
$ws.Range("K24").Value2 = 0
$ws.Range("L24").Value2 = 0
$ws.Range("M24").Value2 = 0
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0
$ws.Range("G25").Value2 = 0
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("M25").Value2 = 0
$ws.Range("D26").Value2 = 23500
$ws.Range("E26").Value2 = 35400
$ws.Range("F26").Value2 = 37000
$ws.Range("G26").Value2 = 32600
$ws.Range("H26").Value2 = 32000
$ws.Range("I26").Value2 = 28500
$ws.Range("J26").Value2 = 43900
$ws.Range("K26").Value2 = 21800
$ws.Range("L26").Value2 = 22500
$ws.Range("M26").Value2 = 35800
$ws.Range("D27").Value2 = 25500
$ws.Range("E27").Value2 = 30400
$ws.Range("F27").Value2 = 31500
$ws.Range("G27").Value2 = 28100
$ws.Range("H27").Value2 = 29300
$ws.Range("I27").Value2 = 25500
$ws.Range("J27").Value2 = 39100
$ws.Range("K27").Value2 = 19000
$ws.Range("L27").Value2 = 19600
$ws.Range("M27").Value2 = 31700
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Value2 = 0
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("L28").Value2 = 0
$ws.Range("M28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("F29").Value2 = 0
$ws.Range("G29").Value2 = 0
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("F30").Value2 = 0
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("M30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 0
$ws.Range("G31").Value2 = 0
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = 0
$ws.Range("D32").Value2 = -900
$ws.Range("E32").Value2 = -1400
$ws.Range("F32").Value2 = -700
$ws.Range("G32").Value2 = -900
$ws.Range("H32").Value2 = -600
$ws.Range("I32").Value2 = -1100
$ws.Range("J32").Value2 = -1000
$ws.Range("K32").Value2 = -400
$ws.Range("L32").Value2 = -600
$ws.Range("M32").Value2 = -12600
$ws.Range("D33").Value2 = 25500
$ws.Range("E33").Value2 = 30400
$ws.Range("F33").Value2 = 31500
$ws.Range("G33").Value2 = 28100
$ws.Range("H33").Value2 = 29300
$ws.Range("I33").Value2 = 25500
$ws.Range("J33").Value2 = 39100
$ws.Range("K33").Value2 = 19000
$ws.Range("L33").Value2 = 19600
$ws.Range("M33").Value2 = 31700
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("F34").Value2 = 0
$ws.Range("G34").Value2 = 0
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = 0
$ws.Range("D35").Value2 = 25500
$ws.Range("E35").Value2 = 30400
$ws.Range("F35").Value2 = 31500
$ws.Range("G35").Value2 = 28100
$ws.Range("H35").Value2 = 29300
$ws.Range("I35").Value2 = 25500
$ws.Range("J35").Value2 = 39100
$ws.Range("K35").Value2 = 19000
$ws.Range("L35").Value2 = 19600
$ws.Range("M35").Value2 = 31700
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("F38").Value2 = 43281
$ws.Range("G38").Value2 = 43190
$ws.Range("H38").Value2 = 43100
$ws.Range("I38").Value2 = 43008
$ws.Range("J38").Value2 = 42916
$ws.Range("K38").Value2 = 42825
$ws.Range("L38").Value2 = 42735
$ws.Range("M38").Value2 = 42643
$ws.Range("D41").Value2 = 146200
$ws.Range("E41").Value2 = 172400
$ws.Range("F41").Value2 = 170400
$ws.Range("G41").Value2 = 183600
$ws.Range("H41").Value2 = 176600
$ws.Range("I41").Value2 = 167700
$ws.Range("J41").Value2 = 173200
$ws.Range("K41").Value2 = 85500
$ws.Range("L41").Value2 = 112900
$ws.Range("M41").Value2 = 158400
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("F42").Value2 = 0
$ws.Range("G42").Value2 = 0
$ws.Range("H42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = 0
$ws.Range("L42").Value2 = 0
$ws.Range("M42").Value2 = 0
$ws.Range("D43").Value2 = 129200
$ws.Range("E43").Value2 = 123800
$ws.Range("F43").Value2 = 118500
$ws.Range("G43").Value2 = 113900
$ws.Range("H43").Value2 = 109000
$ws.Range("I43").Value2 = 105000
$ws.Range("J43").Value2 = 101700
$ws.Range("K43").Value2 = 100100
$ws.Range("L43").Value2 = 95300
$ws.Range("M43").Value2 = 92600
$ws.Range("D44").Value2 = "NA"
$ws.Range("E44").Value2 = "NA"
$ws.Range("F44").Value2 = "NA"
$ws.Range("G44").Value2 = 17600
$ws.Range("H44").Value2 = "NA"
$ws.Range("I44").Value2 = "NA"
$ws.Range("J44").Value2 = "NA"
$ws.Range("K44").Value2 = "NA"
$ws.Range("L44").Value2 = 0
$ws.Range("M44").Value2 = 0
$ws.Range("D45").Value2 = 7800
$ws.Range("E45").Value2 = 12300
$ws.Range("F45").Value2 = 3000
$ws.Range("G45").Value2 = 7900
$ws.Range("H45").Value2 = 9200
$ws.Range("I45").Value2 = 11500
$ws.Range("J45").Value2 = 2600
$ws.Range("K45").Value2 = 7300
$ws.Range("L45").Value2 = 6800
$ws.Range("M45").Value2 = 11300
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("F46").Value2 = 0
$ws.Range("G46").Value2 = 0
$ws.Range("H46").Value2 = 0
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 0
$ws.Range("M46").Value2 = 0
$ws.Range("D47").Value2 = 98800
$ws.Range("E47").Value2 = 97400
$ws.Range("F47").Value2 = 90500
$ws.Range("G47").Value2 = 93500
$ws.Range("H47").Value2 = 98300
$ws.Range("I47").Value2 = 100400
$ws.Range("J47").Value2 = 102600
$ws.Range("K47").Value2 = 138800
$ws.Range("L47").Value2 = 144300
$ws.Range("M47").Value2 = 144900
$ws.Range("D48").Value2 = 7787500
$ws.Range("E48").Value2 = 7787900
$ws.Range("F48").Value2 = 7788300
$ws.Range("G48").Value2 = 7780500
$ws.Range("H48").Value2 = 7817600
$ws.Range("I48").Value2 = 7677400
$ws.Range("J48").Value2 = 7514800
$ws.Range("K48").Value2 = 7176900
$ws.Range("L48").Value2 = 7209500
$ws.Range("M48").Value2 = 7227900
$ws.Range("D49").Value2 = 5200
$ws.Range("E49").Value2 = 5500
$ws.Range("F49").Value2 = 5700
$ws.Range("G49").Value2 = 6000
$ws.Range("H49").Value2 = 6300
$ws.Range("I49").Value2 = 6600
$ws.Range("J49").Value2 = 6700
$ws.Range("K49").Value2 = 6800
$ws.Range("L49").Value2 = 7100
$ws.Range("M49").Value2 = 7000
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("F50").Value2 = 0
$ws.Range("G50").Value2 = 0
$ws.Range("H50").Value2 = 0
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("K50").Value2 = 0
$ws.Range("L50").Value2 = 0
$ws.Range("M50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("F51").Value2 = 0
$ws.Range("G51").Value2 = 0
$ws.Range("H51").Value2 = 0
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("K51").Value2 = 0
$ws.Range("L51").Value2 = 0
$ws.Range("M51").Value2 = 0
$ws.Range("D52").Value2 = 100
$ws.Range("E52").Value2 = 100
$ws.Range("F52").Value2 = 100
$ws.Range("G52").Value2 = 100
$ws.Range("H52").Value2 = 100
$ws.Range("I52").Value2 = 100
$ws.Range("J52").Value2 = 10100
$ws.Range("K52").Value2 = 24100
$ws.Range("L52").Value2 = 100
$ws.Range("M52").Value2 = 200
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("F53").Value2 = 0
$ws.Range("G53").Value2 = 0
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 0
$ws.Range("L53").Value2 = 0
$ws.Range("M53").Value2 = 0
$ws.Range("D54").Value2 = 8261700
$ws.Range("E54").Value2 = 8347800
$ws.Range("F54").Value2 = 8318000
$ws.Range("G54").Value2 = 8323900
$ws.Range("H54").Value2 = 8292600
$ws.Range("I54").Value2 = 8113200
$ws.Range("J54").Value2 = 7955400
$ws.Range("K54").Value2 = 7586700
$ws.Range("L54").Value2 = 7613700
$ws.Range("M54").Value2 = 7645300
$ws.Range("D57").Value2 = 75100
$ws.Range("E57").Value2 = 102400
$ws.Range("F57").Value2 = 75700
$ws.Range("G57").Value2 = 90900
$ws.Range("H57").Value2 = 62700
$ws.Range("I57").Value2 = 80600
$ws.Range("J57").Value2 = 51600
$ws.Range("K57").Value2 = 57100
$ws.Range("L57").Value2 = 36900
$ws.Range("M57").Value2 = 59500
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("F58").Value2 = 0
$ws.Range("G58").Value2 = 0
$ws.Range("H58").Value2 = 0
$ws.Range("I58").Value2 = 0
$ws.Range("J58").Value2 = 0
$ws.Range("K58").Value2 = 0
$ws.Range("L58").Value2 = 0
$ws.Range("M58").Value2 = 0
$ws.Range("D59").Value2 = 105700
$ws.Range("E59").Value2 = 103400
$ws.Range("F59").Value2 = 103600
$ws.Range("G59").Value2 = 103300
$ws.Range("H59").Value2 = 102600
$ws.Range("I59").Value2 = 97500
$ws.Range("J59").Value2 = 96300
$ws.Range("K59").Value2 = 92400
$ws.Range("L59").Value2 = 90400
$ws.Range("M59").Value2 = 89600
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("F60").Value2 = 0
$ws.Range("G60").Value2 = 0
$ws.Range("H60").Value2 = 0
$ws.Range("I60").Value2 = 0
$ws.Range("J60").Value2 = 0
$ws.Range("K60").Value2 = 0
$ws.Range("L60").Value2 = 0
$ws.Range("M60").Value2 = 0
$ws.Range("D61").Value2 = 4134000
$ws.Range("E61").Value2 = 4122300
$ws.Range("F61").Value2 = 4106500
$ws.Range("G61").Value2 = 4098900
$ws.Range("H61").Value2 = 4117400
$ws.Range("I61").Value2 = 4048800
$ws.Range("J61").Value2 = 4314100
$ws.Range("K61").Value2 = 4391400
$ws.Range("L61").Value2 = 4369500
$ws.Range("M61").Value2 = 4401900
$ws.Range("D62").Value2 = 0
$ws.Range("E62").Value2 = 0
$ws.Range("F62").Value2 = 0
$ws.Range("G62").Value2 = 0
$ws.Range("H62").Value2 = 0
$ws.Range("I62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("F63").Value2 = 0
$ws.Range("G63").Value2 = 0
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("M63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("F64").Value2 = 0
$ws.Range("G64").Value2 = 0
$ws.Range("H64").Value2 = 0
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("K64").Value2 = 0
$ws.Range("L64").Value2 = 0
$ws.Range("M64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("F65").Value2 = 0
$ws.Range("G65").Value2 = 0
$ws.Range("H65").Value2 = 0
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = 0
$ws.Range("D66").Value2 = 5859400
$ws.Range("E66").Value2 = 5890700
$ws.Range("F66").Value2 = 5856200
$ws.Range("G66").Value2 = 5865000
$ws.Range("H66").Value2 = 5855100
$ws.Range("I66").Value2 = 5677300
$ws.Range("J66").Value2 = 5828300
$ws.Range("K66").Value2 = 5722200
$ws.Range("L66").Value2 = 5692600
$ws.Range("M66").Value2 = 5753100
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("F68").Value2 = 0
$ws.Range("G68").Value2 = 0
$ws.Range("H68").Value2 = 0
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("K68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("M68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 0
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("M69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 0
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 0
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("M70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("F71").Value2 = 0
$ws.Range("G71").Value2 = 0
$ws.Range("H71").Value2 = 0
$ws.Range("I71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("K71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("M71").Value2 = 0
$ws.Range("D72").Value2 = -935600
$ws.Range("E72").Value2 = -917000
$ws.Range("F72").Value2 = -905100
$ws.Range("G72").Value2 = -894300
$ws.Range("H72").Value2 = -879800
$ws.Range("I72").Value2 = -867000
$ws.Range("J72").Value2 = -853600
$ws.Range("K72").Value2 = -836900
$ws.Range("L72").Value2 = -820700
$ws.Range("M72").Value2 = -805500
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("F73").Value2 = 0
$ws.Range("G73").Value2 = 0
$ws.Range("H73").Value2 = 0
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("M73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("F74").Value2 = 0
$ws.Range("G74").Value2 = 0
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 0
$ws.Range("G75").Value2 = 0
$ws.Range("H75").Value2 = 0
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("K75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("M75").Value2 = 0
$ws.Range("D76").Value2 = 2402300
$ws.Range("E76").Value2 = 2457100
$ws.Range("F76").Value2 = 2461700
$ws.Range("G76").Value2 = 2458900
$ws.Range("H76").Value2 = 2437500
$ws.Range("I76").Value2 = 2435900
$ws.Range("J76").Value2 = 2127100
$ws.Range("K76").Value2 = 1864500
$ws.Range("L76").Value2 = 1921100
$ws.Range("M76").Value2 = 1892200
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("F77").Value2 = 0
$ws.Range("G77").Value2 = 0
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("F80").Value2 = 43281
$ws.Range("G80").Value2 = 43190
$ws.Range("H80").Value2 = 43100
$ws.Range("I80").Value2 = 43008
$ws.Range("J80").Value2 = 42916
$ws.Range("K80").Value2 = 42825
$ws.Range("L80").Value2 = 42735
$ws.Range("M80").Value2 = 42643
$ws.Range("D81").Value2 = 25500
$ws.Range("E81").Value2 = 30400
$ws.Range("F81").Value2 = 31500
$ws.Range("G81").Value2 = 28100
$ws.Range("H81").Value2 = 29300
$ws.Range("I81").Value2 = 25500
$ws.Range("J81").Value2 = 39100
$ws.Range("K81").Value2 = 19000
$ws.Range("L81").Value2 = 19600
$ws.Range("M81").Value2 = 31700
$ws.Range("D83").Value2 = 89900
$ws.Range("E83").Value2 = 74100
$ws.Range("F83").Value2 = 73400
$ws.Range("G83").Value2 = 72500
$ws.Range("H83").Value2 = 70600
$ws.Range("I83").Value2 = 70000
$ws.Range("J83").Value2 = 136200
$ws.Range("K83").Value2 = 67400
$ws.Range("L83").Value2 = 67000
$ws.Range("M83").Value2 = 63800
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("F84").Value2 = 0
$ws.Range("G84").Value2 = 0
$ws.Range("H84").Value2 = 0
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 0
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("F85").Value2 = 0
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("L85").Value2 = 0
$ws.Range("M85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("F86").Value2 = 0
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 0
$ws.Range("M86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("F87").Value2 = 0
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 0
$ws.Range("I87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("K87").Value2 = 0
$ws.Range("L87").Value2 = 0
$ws.Range("M87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("F88").Value2 = 0
$ws.Range("G88").Value2 = 0
$ws.Range("H88").Value2 = 0
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("K88").Value2 = 0
$ws.Range("L88").Value2 = 0
$ws.Range("M88").Value2 = 0
$ws.Range("D89").Value2 = 101100
$ws.Range("E89").Value2 = 116000
$ws.Range("F89").Value2 = 98300
$ws.Range("G89").Value2 = 117600
$ws.Range("H89").Value2 = 94000
$ws.Range("I89").Value2 = 110100
$ws.Range("J89").Value2 = 198600
$ws.Range("K89").Value2 = 107600
$ws.Range("L89").Value2 = 73800
$ws.Range("M89").Value2 = 98900
$ws.Range("D91").Value2 = -20000
$ws.Range("E91").Value2 = -22000
$ws.Range("F91").Value2 = -15500
$ws.Range("G91").Value2 = -11000
$ws.Range("H91").Value2 = -20900
$ws.Range("I91").Value2 = -32100
$ws.Range("J91").Value2 = -10000
$ws.Range("K91").Value2 = -59200
$ws.Range("L91").Value2 = -43100
$ws.Range("M91").Value2 = -392500
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 0
$ws.Range("G92").Value2 = 0
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("M92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("F93").Value2 = 0
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("M93").Value2 = 0
$ws.Range("D94").Value2 = -81600
$ws.Range("E94").Value2 = -72200
$ws.Range("F94").Value2 = -61500
$ws.Range("G94").Value2 = -34300
$ws.Range("H94").Value2 = -54400
$ws.Range("I94").Value2 = -213500
$ws.Range("J94").Value2 = -401600
$ws.Range("K94").Value2 = -57800
$ws.Range("L94").Value2 = -40900
$ws.Range("M94").Value2 = -282200
$ws.Range("D96").Value2 = -42500
$ws.Range("E96").Value2 = -42500
$ws.Range("F96").Value2 = -42500
$ws.Range("G96").Value2 = -42400
$ws.Range("H96").Value2 = -39000
$ws.Range("I96").Value2 = -37000
$ws.Range("J96").Value2 = -70100
$ws.Range("K96").Value2 = -34900
$ws.Range("L96").Value2 = -33200
$ws.Range("M96").Value2 = -32800
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("F97").Value2 = 0
$ws.Range("G97").Value2 = 0
$ws.Range("H97").Value2 = 0
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 0
$ws.Range("L97").Value2 = 0
$ws.Range("M97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("F98").Value2 = 0
$ws.Range("G98").Value2 = 0
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("M98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("F99").Value2 = 0
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = 0
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("L99").Value2 = 0
$ws.Range("M99").Value2 = 0
$ws.Range("D100").Value2 = -45700
$ws.Range("E100").Value2 = -41800
$ws.Range("F100").Value2 = -49900
$ws.Range("G100").Value2 = -76400
$ws.Range("H100").Value2 = -30700
$ws.Range("I100").Value2 = 98000
$ws.Range("J100").Value2 = 263300
$ws.Range("K100").Value2 = -77200
$ws.Range("L100").Value2 = -78300
$ws.Range("M100").Value2 = 264600
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("F101").Value2 = 0
$ws.Range("G101").Value2 = 0
$ws.Range("H101").Value2 = 0
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 0
$ws.Range("K101").Value2 = 0
$ws.Range("L101").Value2 = 0
$ws.Range("M101").Value2 = 0
$ws.Range("D102").Value2 = -26100
$ws.Range("E102").Value2 = 2000
$ws.Range("F102").Value2 = -13200
$ws.Range("G102").Value2 = 6900
$ws.Range("H102").Value2 = 8900
$ws.Range("I102").Value2 = -5400
$ws.Range("J102").Value2 = 60200
$ws.Range("K102").Value2 = -27400
$ws.Range("L102").Value2 = -45500
$ws.Range("M102").Value2 = 81200

Write-Output "Workbook updated: inserted 2 new quarterly columns and refreshed financial figures."
